$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Popis dokumentacije ")

# Row 28 was blank; fill it in with the new PB_23 / UC030 (PDV registar) entry.
# Column A ("23.") reads like a number to the auto-type-detection, so enter it
# as a formula first and convert the result to a plain value - this keeps it
# as text (matching the other "NN." cells in the column) without disturbing
# the cell's existing style.
$ws.Range("A28").Formula = '="23."'
$ws.Range("A28").Copy()
$ws.Range("A28").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B28").Value = "PB_23"
$ws.Range("E28").Value = "08.05.2014."
$ws.Range("C28").Value = "Specifikacija UC030 Upravljanje funkcionalnostima PDV registra"
$ws.Range("D28").Value = "v 1.0"
$ws.Range("F28").Value = "Specifikacija slučajeva korištenja"

# The long, wrapped description makes the new row two lines tall
$ws.Rows(28).RowHeight = 30

# Row 26: shorten the UC102 document title (drop "slučajeva korištenja " from the name)
$ws.Range("C26").Value = "Specifikacija UC102 Upravljanjem funkcionalnostima aplikacije"

# Move the active selection, matching the author's final cursor position
$ws.Range("C33").Select() | Out-Null
